# Generate Report for Handback
#
# This run re-generates the handback-status report: the in-flight file
# (4c08c1d8-...) finished and is replaced by a freshly-generated file
# (2db4db43-...), and a second file (dccdc1f8-...) was also handed back,
# landing as a new row 3 on every sheet / table.

$wb = $excel.ActiveWorkbook

$oldGuid = "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc"
$newGuid1 = "2db4db43-ed49-4db6-94a7-647e8cb93e42"
$newGuid2 = "dccdc1f8-1fc0-4ee0-ac3a-7fe9156b3bc5"

$zhHash1 = "7c2217b3072ade006b225d85acc5ed30726c04e1"
$deHash1 = "7c2217b3072ade006b225d85acc5ed30726c04e1"
$zhHash2 = "5e35570a73f1e5e1d75101395355b124d563c1cf"
$deHash2 = "5e35570a73f1e5e1d75101395355b124d563c1cf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = "e2e\$newGuid1.md"
$ws.Range("C2").Value = ".md"
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Range("G2").Value = "2016-08-30 11:10:32"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "e2e\$newGuid2.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-30 11:10:32"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$newGuid1.md", "", "", "e2e\$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$newGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "$newGuid1.$zhHash1.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-30 11:10:27"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$zhHash1.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-30 11:10:55"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = "'"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$newGuid2.$zhHash2.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-30 11:10:27"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid2.$zhHash2.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-30 11:10:55"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/449786648024e3023869a4bb1bfb790c8b6f90d2/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/449786648024e3023869a4bb1bfb790c8b6f90d2/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

$ws.Range("A2").Style = "Hyperlink"
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "$newGuid1.$deHash1.de-de.xlf"
$ws.Range("H2").Value = "2016-08-30 11:10:32"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$deHash1.de-de.xlf"
$ws.Range("K2").Value = "2016-08-30 11:11:07"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$newGuid2.$deHash2.de-de.xlf"
$ws.Range("H3").Value = "2016-08-30 11:10:32"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid2.$deHash2.de-de.xlf"
$ws.Range("K3").Value = "2016-08-30 11:11:07"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8389c72fbaabd0b088a90bcf0be01e37d832d832/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8389c72fbaabd0b088a90bcf0be01e37d832d832/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

$ws.Range("A2").Style = "Hyperlink"
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"

Write-Host "Report regenerated: handback rows for $newGuid1 and $newGuid2 written across Overview/zh-cn/de-de."
